$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 285, shifting rows 285:347 down to 286:347
$ws.Rows.Item(285).Insert()

$ws.Range("A285").Value = 3
$ws.Range("B285").Value = "Femacal de La Calera"
$ws.Range("C285").Value = "Coquimbo"
$ws.Range("D285").Value = 44711
$ws.Range("E285").Value = 5
$ws.Range("F285").Value = 100112012
$ws.Range("G285").Value = "Espinaca"
$ws.Range("H285").Value = "Sin especificar"
$ws.Range("I285").Value = "Primera"
$ws.Range("J285").Value = 230
$ws.Range("K285").Value = 3500
$ws.Range("L285").Value = 4000
$ws.Range("M285").Value = 3739
$ws.Range("N285").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O285").Value = "Provincia de Quillota"
$ws.Range("P285").Value = 1246
$ws.Range("Q285").Value = 3
$ws.Range("R285").Value = "Hortaliza"
